$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Prepare new rows 6-12 by copying formatting (styles) from existing
# rows that already have the right per-column styles, then fill values.
# xlPasteFormats = -4122
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

# Rows 6-9 follow the same per-column style pattern as row 5 (B:M)
$ws.Range("B5:M5").Copy() | Out-Null
$ws.Range("B6:M9").PasteSpecial($xlPasteFormats) | Out-Null

# Rows 10-11 follow the same per-column style pattern as row 5 but without
# values/styles in G:J (those columns stay completely blank/unstyled)
$ws.Range("B5:F5").Copy() | Out-Null
$ws.Range("B10:F11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K5:M5").Copy() | Out-Null
$ws.Range("K10:M11").PasteSpecial($xlPasteFormats) | Out-Null

# Row 12 follows the same per-column style pattern as row 5 across B:M
$ws.Range("B5:M5").Copy() | Out-Null
$ws.Range("B12:M12").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Fill in the brand-new text values first (in the same order the rows
# were authored), then backfill the cells that reuse existing text.
# ---------------------------------------------------------------------

# Row 5 (K/L/M get new content)
$ws.Range("K5").Value = "Constraint Margin Rate based on  Add Tax"
$ws.Range("L5").Value = "marginRateBasedOnFreightAndAddTax"
$ws.Range("M5").Value = "rule"

# Row 6 (new)
$ws.Range("K6").Value = "Test_SR"
$ws.Range("L6").Value = "test_SR"
$ws.Range("M6").Value = "library"

# Row 9 (new)
$ws.Range("K9").Value = "Tax Exempt Reason"
$ws.Range("L9").Value = "taxExemptReason"

# Row 2 (K/L/M get new content)
$ws.Range("K2").Value = "User Import"
$ws.Range("L2").Value = "oRCL_SFDC_UserImport"
$ws.Range("M2").Value = "integration"

# Row 2 col A gets new content
$ws.Range("A2").Value = "Demo28AugT1"

# ---------------------------------------------------------------------
# Row 2 remaining updates
# ---------------------------------------------------------------------
$ws.Range("J2").ClearContents()

# ---------------------------------------------------------------------
# Row 3 updates
# ---------------------------------------------------------------------
$ws.Range("K3").Value = "API_Save"
$ws.Range("L3").Value = "aPI_Save_t"
$ws.Range("M3").Value = "action"

# ---------------------------------------------------------------------
# Row 4 updates
# ---------------------------------------------------------------------
$ws.Range("H4").Value = "Transaction"
$ws.Range("I4").Value = "transaction"
$ws.Range("K4").Value = "Last Priced"
$ws.Range("L4").Value = "lastPricedDate_t"
$ws.Range("M4").Value = "attribute"

# ---------------------------------------------------------------------
# Row 5 remaining updates
# ---------------------------------------------------------------------
$ws.Range("H5").Value = "Transaction"
$ws.Range("I5").Value = "transaction"

# ---------------------------------------------------------------------
# Row 6 remaining (new)
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "Commerce"
$ws.Range("C6").Value = "COMMERCE"
$ws.Range("D6").Value = "Paramount Quote to Order"
$ws.Range("E6").Value = "oraclecpqo_bmClone_2"
$ws.Range("F6").Value = "process"
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = "Transaction"
$ws.Range("I6").Value = "transaction"
$ws.Range("J6").Value = "document"

# ---------------------------------------------------------------------
# Row 7 (new)
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "Commerce"
$ws.Range("C7").Value = "COMMERCE"
$ws.Range("D7").Value = "Paramount Quote to Order"
$ws.Range("E7").Value = "oraclecpqo_bmClone_2"
$ws.Range("F7").Value = "process"
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = "Transaction Line"
$ws.Range("I7").Value = "transactionLine"
$ws.Range("J7").Value = "document"
$ws.Range("K7").Value = "Back"
$ws.Range("L7").Value = "back_l"
$ws.Range("M7").Value = "action"

# ---------------------------------------------------------------------
# Row 8 (new)
# ---------------------------------------------------------------------
$ws.Range("B8").Value = "Commerce"
$ws.Range("C8").Value = "COMMERCE"
$ws.Range("D8").Value = "Paramount Quote to Order"
$ws.Range("E8").Value = "oraclecpqo_bmClone_2"
$ws.Range("F8").Value = "process"
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = "Transaction Line"
$ws.Range("I8").Value = "transactionLine"
$ws.Range("J8").Value = "document"
$ws.Range("K8").Value = "Document Number"
$ws.Range("L8").Value = "documentNumber_l"
$ws.Range("M8").Value = "attribute"

# ---------------------------------------------------------------------
# Row 9 remaining (new)
# ---------------------------------------------------------------------
$ws.Range("B9").Value = "Commerce"
$ws.Range("C9").Value = "COMMERCE"
$ws.Range("D9").Value = "Paramount Quote to Order"
$ws.Range("E9").Value = "oraclecpqo_bmClone_2"
$ws.Range("F9").Value = "process"
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = "Transaction Line"
$ws.Range("I9").Value = "transactionLine"
$ws.Range("J9").Value = "document"
$ws.Range("M9").Value = "rule"

# ---------------------------------------------------------------------
# Row 10 (new) - Document Designer / Field Profile Sheet
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Document Designer"
$ws.Range("C10").Value = "DOCUMENT_DESIGNER"
$ws.Range("D10").Value = "Paramount Quote to Order"
$ws.Range("E10").Value = "oraclecpqo_bmClone_2"
$ws.Range("F10").Value = "_set"
$ws.Range("K10").Value = "Field Profile Sheet - English"
$ws.Range("L10").Value = "Field Profile Sheet - English"
$ws.Range("M10").Value = "doc_designer"

# ---------------------------------------------------------------------
# Row 11 (new) - Document Designer / Job Profile Sheet
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "Document Designer"
$ws.Range("C11").Value = "DOCUMENT_DESIGNER"
$ws.Range("D11").Value = "Paramount Quote to Order"
$ws.Range("E11").Value = "oraclecpqo_bmClone_2"
$ws.Range("F11").Value = "_set"
$ws.Range("K11").Value = "Job Profile Sheet - English"
$ws.Range("L11").Value = "Job Profile Sheet - English"
$ws.Range("M11").Value = "doc_designer"

# ---------------------------------------------------------------------
# Row 12 (new) - Email Designer / Final Approval Notification
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "Email Designer"
$ws.Range("C12").Value = "EMAIL_DESIGNER"
$ws.Range("D12").Value = "Paramount Quote to Order"
$ws.Range("E12").Value = "oraclecpqo_bmClone_2"
$ws.Range("F12").Value = "_set"
$ws.Range("G12").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("K12").Value = "Final Approval Notification For DOA - English"
$ws.Range("L12").Value = "Final Approval Notification For DOA - English"
$ws.Range("M12").Value = "email_designer"

# ---------------------------------------------------------------------
# Sheet view: zoom to 85% and select H25
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("H25").Select() | Out-Null
